# Auto-generated script to update cryptos.xlsx per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.999.68'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '3.257.98'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = "'599.22"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'137.25"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.260.15'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = "'0.509"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').Value = "'0.146"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('D11').Value = "'5.42"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').Value = "'0.459"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').Value = "'0.0000241"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.83%  '
$ws.Range('D14').Value = "'33.85"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').Value = '3.808.51'
$ws.Range('E15').Value = '  +0.48%  '
$ws.Range('E16').Value = '  +1.28%  '
$ws.Range('D17').Value = '3.269.23'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '63.100.19'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').Value = "'6.70"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.50%  '
$ws.Range('D20').Value = "'469.91"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.93%  '
$ws.Range('D21').Value = "'13.75"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.65%  '
$ws.Range('D22').Value = "'0.720"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('D23').Value = "'7.81"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.12%  '
$ws.Range('D24').Value = "'13.52"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('D25').Value = "'83.80"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').Value = "'7.04"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.49%  '
$ws.Range('D30').Value = "'7.89"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.46%  '
$ws.Range('D31').Value = "'2.10"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.99%  '
$ws.Range('D32').Value = "'27.93"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('D33').Value = "'0.103"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.52%  '
$ws.Range('D34').Value = "'2.46"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.96%  '
$ws.Range('D35').Value = "'1.07"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.50%  '
$ws.Range('D36').Value = "'5.87"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').Value = "'51.73"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.01%  '
$ws.Range('D38').Value = '0.0₃0715'
$ws.Range('E38').Value = '  -0.63%  '
$ws.Range('D39').Value = "'0.0393"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('D40').Value = '3.068.86'
$ws.Range('E40').Value = '  +2.36%  '
$ws.Range('D41').Value = "'421.04"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.57%  '
$ws.Range('D42').Value = "'0.116"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.04%  '
$ws.Range('D43').Value = "'8.17"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.63%  '
$ws.Range('D44').Value = "'2.63"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.04%  '
$ws.Range('D45').Value = "'0.256"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.52%  '
$ws.Range('D46').Value = "'2.15"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.27%  '
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = "'127.05"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.47%  '
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').Value = "'35.69"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.46%  '
$ws.Range('D50').Value = "'25.68"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('E51').Value = '  -2.12%  '
